# Atualização do banco e dos requisitos
#
# The "Perfil" backlog item's description cell (D9 in the Tabela1 table on
# Planilha1) was reworded to talk about posts + interactions instead of
# posts + bio + follower count.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("D9").Value = "Página que irá conter os posts usuário e as interações nele contidas."

# Mirror the author's final selection: they had just finished editing D9.
$ws.Range("D9").Select()
